# Update "想去人数" (want-to-go count) values in column F across sheets,
# matching the commit "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 366
$ws1.Range("F4").Value = 276
$ws1.Range("F5").Value = 4106
$ws1.Range("F6").Value = 38
$ws1.Range("F7").Value = 452

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 6

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 366
$ws4.Range("F4").Value = 276
$ws4.Range("F5").Value = 4106
$ws4.Range("F7").Value = 6
$ws4.Range("F8").Value = 38
$ws4.Range("F9").Value = 452
